$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (the existing "Sections" row), shifting
# everything below it down by one. This mirrors splitting the former
# "Wind" block into its own "wind" section (Ground_rou / B) ahead of
# "Sections".
$ws.Rows.Item(5).Insert()

# New row 5: Ground_rou / B, styled like the other parameter rows
# (B3/B4 use style index 1 -> centered alignment). Write B5 first so the
# shared-string table gets "B" before "Ground_rou", matching the target.
$ws.Range("B5").Value = "B"
$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("A5").Value = "Ground_rou"

# Keep the active selection pointing at the newly inserted row, as in
# the target workbook.
$ws.Range("A5").Select()
